$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "Yearly" sheet: October's Taxable Account dividend was corrected from
# 17.55 to 35.24. Everything else on this sheet (G12, D15, G15) is a
# formula that recalculates automatically from this single input edit.
# ---------------------------------------------------------------------
$wsYearly = $wb.Worksheets.Item("Yearly")
$wsYearly.Range("D12").Value = 35.24

# Leftover selection state on the Yearly sheet after the edit.
$wsYearly.Range("I13").Select()

# ---------------------------------------------------------------------
# "All Time" sheet: the 2016 Taxable Account total (F7) mirrors the
# Yearly sheet's D15 total but is stored as a plain value here, so it
# has to be updated explicitly to stay in sync. I7/F46/I46 are formulas
# and recalc automatically.
# ---------------------------------------------------------------------
$wsAllTime = $wb.Worksheets.Item("All Time")
$wsAllTime.Range("F7").Value = 608.14

# Move to the "All Time" sheet (it ends up the active/visible tab) and
# restore its final scroll position + selection.
$wsAllTime.Activate()
$excel.ActiveWindow.ScrollRow = 43
$excel.ActiveWindow.ScrollColumn = 1
$wsAllTime.Range("J52").Select()
